$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 39675
$ws.Range("E2").Value = -8677
$ws.Range("F2").Value = -8677
$ws.Range("G2").Value = -8964
$ws.Range("H2").Value = -6793
$ws.Range("I2").Value = -6328
$ws.Range("J2").Value = -464
$ws.Range("K2").Value = 100637
$ws.Range("L2").Value = 81371
$ws.Range("M2").Value = 19266
$ws.Range("N2").Value = 17973
$ws.Range("O2").Value = 1293
$ws.Range("P2").Value = 1000
$ws.Range("Q2").Value = 4058
$ws.Range("R2").Value = 1811
$ws.Range("S2").Value = -1139
$ws.Range("T2").Value = 627
$ws.Range("U2").Value = 3430
$ws.Range("V2").Value = 53930
$ws.Range("W2").Value = -21.87
$ws.Range("X2").Value = -17.12
$ws.Range("Y2").Value = -25.88
$ws.Range("Z2").Value = -6.69
$ws.Range("AA2").Value = 422.36
$ws.Range("AB2").Value = 1302.19
$ws.Range("AC2").Value = -15844
$ws.Range("AE2").Value = 45299
$ws.Range("AF2").Value = 0.77
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 39942149
$ws.Range("AD2").Value = $null

# Row 3
$ws.Range("D3").Value = 46524
$ws.Range("E3").Value = 666
$ws.Range("F3").Value = 666
$ws.Range("G3").Value = 498
$ws.Range("H3").Value = 257
$ws.Range("I3").Value = 382
$ws.Range("J3").Value = -125
$ws.Range("K3").Value = 94380
$ws.Range("L3").Value = 76412
$ws.Range("M3").Value = 17968
$ws.Range("N3").Value = 16750
$ws.Range("O3").Value = 1218
$ws.Range("P3").Value = 1000
$ws.Range("Q3").Value = -12126
$ws.Range("R3").Value = -657
$ws.Range("S3").Value = 12184
$ws.Range("T3").Value = 367
$ws.Range("U3").Value = -12493
$ws.Range("V3").Value = 53997
$ws.Range("W3").Value = 1.43
$ws.Range("X3").Value = 0.55
$ws.Range("Y3").Value = 2.2
$ws.Range("Z3").Value = 0.26
$ws.Range("AA3").Value = 425.27
$ws.Range("AB3").Value = 1325.38
$ws.Range("AC3").Value = 956
$ws.Range("AE3").Value = 41996
$ws.Range("AF3").Value = 0.62
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 39942149
$ws.Range("AD3").Value = $null

# Row 4
$ws.Range("D4").Value = 34465
$ws.Range("E4").Value = 1911
$ws.Range("F4").Value = 2075
$ws.Range("G4").Value = 1910
$ws.Range("H4").Value = 396
$ws.Range("I4").Value = 346
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 92263
$ws.Range("L4").Value = 69679
$ws.Range("M4").Value = 22584
$ws.Range("N4").Value = 21425
$ws.Range("O4").Value = 1159
$ws.Range("P4").Value = 1000
$ws.Range("Q4").Value = 5297
$ws.Range("R4").Value = 1076
$ws.Range("S4").Value = -4899
$ws.Range("T4").Value = 302
$ws.Range("U4").Value = 4994
$ws.Range("V4").Value = 53533
$ws.Range("W4").Value = 5.54
$ws.Range("X4").Value = 1.15
$ws.Range("Y4").Value = 1.81
$ws.Range("Z4").Value = 0.42
$ws.Range("AA4").Value = 308.54
$ws.Range("AB4").Value = 1362.19
$ws.Range("AC4").Value = 866
$ws.Range("AE4").Value = 53717
$ws.Range("AF4").Value = 0.63
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 39942149
$ws.Range("AD4").Value = $null

# Row 5
$ws.Range("D5").Value = 24534
$ws.Range("E5").Value = 1079
$ws.Range("F5").Value = 1079
$ws.Range("G5").Value = 7333
$ws.Range("H5").Value = 4452
$ws.Range("I5").Value = 4369
$ws.Range("J5").Value = 82
$ws.Range("K5").Value = 83194
$ws.Range("L5").Value = 59803
$ws.Range("M5").Value = 23391
$ws.Range("N5").Value = 22168
$ws.Range("O5").Value = 1224
$ws.Range("P5").Value = 1000
$ws.Range("Q5").Value = -6793
$ws.Range("R5").Value = 7414
$ws.Range("S5").Value = -2050
$ws.Range("T5").Value = 654
$ws.Range("U5").Value = -7448
$ws.Range("V5").Value = 4624
$ws.Range("W5").Value = 4.4
$ws.Range("X5").Value = 18.14
$ws.Range("Y5").Value = 20.05
$ws.Range("Z5").Value = 5.07
$ws.Range("AA5").Value = 255.66
$ws.Range("AB5").Value = 1801.36
$ws.Range("AC5").Value = 10939
$ws.Range("AD5").Value = 3.61
$ws.Range("AE5").Value = 55580
$ws.Range("AF5").Value = 0.71
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 39942149

# Row 6
$ws.Range("D6").Value = 24030
$ws.Range("E6").Value = 709
$ws.Range("F6").Value = 709
$ws.Range("G6").Value = 1175
$ws.Range("H6").Value = 1207
$ws.Range("I6").Value = 850
$ws.Range("K6").Value = 34637
$ws.Range("L6").Value = 11466
$ws.Range("M6").Value = 23172
$ws.Range("N6").Value = 22896
$ws.Range("P6").Value = 1997
$ws.Range("Q6").Value = 5185
$ws.Range("R6").Value = -107
$ws.Range("S6").Value = -7977
$ws.Range("T6").Value = 4771
$ws.Range("U6").Value = 414
$ws.Range("V6").Value = 2830
$ws.Range("W6").Value = 2.95
$ws.Range("X6").Value = 5.02
$ws.Range("Y6").Value = 3.77
$ws.Range("Z6").Value = 2.05
$ws.Range("AA6").Value = 49.48
$ws.Range("AB6").Value = 950.86
$ws.Range("AC6").Value = 2128
$ws.Range("AD6").Value = 28.19
$ws.Range("AE6").Value = 57405
$ws.Range("AF6").Value = 1.05
$ws.Range("AG6").Value = 700
$ws.Range("AH6").Value = 1.17
$ws.Range("AI6").Value = 32.84
$ws.Range("AJ6").Value = 39942149

# Row 7
$ws.Range("D7").Value = 30329
$ws.Range("E7").Value = 1090
$ws.Range("G7").Value = 1147
$ws.Range("H7").Value = 911
$ws.Range("I7").Value = 807
$ws.Range("K7").Value = 35959
$ws.Range("L7").Value = 12213
$ws.Range("M7").Value = 23717
$ws.Range("N7").Value = 23332
$ws.Range("P7").Value = 1999
$ws.Range("Q7").Value = 1930
$ws.Range("R7").Value = 1519
$ws.Range("S7").Value = -2902
$ws.Range("T7").Value = 636
$ws.Range("U7").Value = 331
$ws.Range("W7").Value = 3.59
$ws.Range("X7").Value = 3
$ws.Range("Y7").Value = 3.49
$ws.Range("Z7").Value = 2.58
$ws.Range("AA7").Value = 51.49
$ws.Range("AC7").Value = 2020
$ws.Range("AD7").Value = 20.52
$ws.Range("AE7").Value = 58498
$ws.Range("AF7").Value = 0.71
$ws.Range("AG7").Value = 636
$ws.Range("AH7").Value = 1.54
$ws.Range("AI7").Value = 31.51

# Row 8
$ws.Range("D8").Value = 30788
$ws.Range("E8").Value = 1240
$ws.Range("G8").Value = 1346
$ws.Range("H8").Value = 1048
$ws.Range("I8").Value = 947
$ws.Range("K8").Value = 37123
$ws.Range("L8").Value = 12618
$ws.Range("M8").Value = 24451
$ws.Range("N8").Value = 23983
$ws.Range("P8").Value = 1999
$ws.Range("Q8").Value = 1559
$ws.Range("R8").Value = -228
$ws.Range("S8").Value = -2436
$ws.Range("T8").Value = 550
$ws.Range("U8").Value = 814
$ws.Range("W8").Value = 4.03
$ws.Range("X8").Value = 3.4
$ws.Range("Y8").Value = 4
$ws.Range("Z8").Value = 2.87
$ws.Range("AA8").Value = 51.6
$ws.Range("AC8").Value = 2371
$ws.Range("AD8").Value = 17.49
$ws.Range("AE8").Value = 60132
$ws.Range("AF8").Value = 0.6899999999999999
$ws.Range("AG8").Value = 668
$ws.Range("AH8").Value = 1.61
$ws.Range("AI8").Value = 28.19

# Row 9
$ws.Range("D9").Value = 33349
$ws.Range("E9").Value = 1644
$ws.Range("G9").Value = 1778
$ws.Range("H9").Value = 1382
$ws.Range("I9").Value = 1239
$ws.Range("K9").Value = 38974
$ws.Range("L9").Value = 13382
$ws.Range("M9").Value = 25507
$ws.Range("N9").Value = 24961
$ws.Range("P9").Value = 1999
$ws.Range("Q9").Value = 1713
$ws.Range("R9").Value = -248
$ws.Range("S9").Value = -2396
$ws.Range("T9").Value = 605
$ws.Range("U9").Value = 836
$ws.Range("W9").Value = 4.93
$ws.Range("X9").Value = 4.15
$ws.Range("Y9").Value = 5.06
$ws.Range("Z9").Value = 3.63
$ws.Range("AA9").Value = 52.46
$ws.Range("AC9").Value = 3101
$ws.Range("AD9").Value = 13.37
$ws.Range("AE9").Value = 62584
$ws.Range("AF9").Value = 0.66
$ws.Range("AG9").Value = 700
$ws.Range("AH9").Value = 22.57
